$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.04897152071431162
$ws.Range("H2").Value = 15.83035461096923
$ws.Range("I2").Value = 0.3745941685984026
$ws.Range("G3").Value = 0.06161644760434993
$ws.Range("H3").Value = 26.52157245051975
$ws.Range("G4").Value = -0.006093676716984563
$ws.Range("H4").Value = -420.063603451725
$ws.Range("G5").Value = 0.002707114828711582
$ws.Range("H5").Value = -28.69165986592789
$ws.Range("G6").Value = 0.02123993786832681
$ws.Range("H6").Value = -38.73609602108809
$ws.Range("G7").Value = 0.07393751212338077
$ws.Range("H7").Value = 39.00916837117111
$ws.Range("G8").Value = -0.02444832773774197
$ws.Range("H8").Value = -29.90074906944205
$ws.Range("G9").Value = 0.02894502473046889
$ws.Range("H9").Value = 234.4535137598322
$ws.Range("G10").Value = -0.08431391100214167
$ws.Range("H10").Value = -15.97295516839898
$ws.Range("G11").Value = -0.08399990927501251
$ws.Range("H11").Value = 8.706707501493383
$ws.Range("G12").Value = -0.228404069243767
$ws.Range("H12").Value = 6.565870344929901
$ws.Range("G13").Value = -0.2705652763176778
$ws.Range("H13").Value = 1.545757002950137
$ws.Range("G14").Value = -0.07851438935273379
$ws.Range("H14").Value = -111.6450914909133
$ws.Range("G15").Value = -0.004228536505660871
$ws.Range("H15").Value = 87.83902444403785
$ws.Range("G16").Value = 0.112293905996712
$ws.Range("H16").Value = -10.38581892984485
$ws.Range("G17").Value = 0.131772785580724
$ws.Range("H17").Value = -6.044688706481937
$ws.Range("G18").Value = 0.1275524207718969
$ws.Range("H18").Value = 2.260591026879167
$ws.Range("G19").Value = 0.1233227248906218
$ws.Range("H19").Value = -7.424992608679556
$ws.Range("G20").Value = 0.04343994456343899
$ws.Range("H20").Value = 26.51386871957572
$ws.Range("G21").Value = 0.04417070683778238
$ws.Range("H21").Value = -23.89591420688458
$ws.Range("G22").Value = -0.06929198704654332
$ws.Range("H22").Value = 13.22011718107524
$ws.Range("G23").Value = -0.05384497281633821
$ws.Range("H23").Value = 13.91200262641725
$ws.Range("G24").Value = 0.1081806314496519
$ws.Range("H24").Value = -8.410379541679704
$ws.Range("G25").Value = 0.1319581053077023
$ws.Range("H25").Value = 4.587218439099161
$ws.Range("G26").Value = 0.06097260409071388
$ws.Range("H26").Value = 22.67115595634789
$ws.Range("G27").Value = 0.08627964584715897
$ws.Range("H27").Value = -0.4579374277349867
$ws.Range("G28").Value = -0.09417325083052169
$ws.Range("H28").Value = -48.0850344366213
$ws.Range("G29").Value = -0.1087870678671063
$ws.Range("H29").Value = -52.85046068234131
$ws.Range("G30").Value = 0.05770942419999141
$ws.Range("H30").Value = -9.415914257693657
$ws.Range("G31").Value = 0.06150796584682331
$ws.Range("H31").Value = 1.531029250733402
$ws.Range("G32").Value = 0.07221310413608721
$ws.Range("H32").Value = -26.51111729106108
$ws.Range("G33").Value = 0.1084390441692693
$ws.Range("H33").Value = 31.78517683584134
$ws.Range("G34").Value = -0.003873312572378081
$ws.Range("H34").Value = -114.8657857126932
$ws.Range("G35").Value = -0.006222015626087195
$ws.Range("H35").Value = 44.49638546615682
$ws.Range("G36").Value = -0.007127508990366663
$ws.Range("H36").Value = -1424.306059271223
$ws.Range("G37").Value = -0.002929029570126664
$ws.Range("H37").Value = 76.66904765975055
$ws.Range("G38").Value = 0.1088286692939782
$ws.Range("H38").Value = 1.465077285340419
$ws.Range("G39").Value = 0.09638141018461149
$ws.Range("H39").Value = 12.51258477727404
$ws.Range("G40").Value = 0.02254409658196638
$ws.Range("H40").Value = 658.9947401429642
$ws.Range("G41").Value = 0.02142627828945714
$ws.Range("H41").Value = 42.8728450227983
$ws.Range("G42").Value = 0.1023554189459312
$ws.Range("H42").Value = 1.407134942014298
$ws.Range("G43").Value = 0.136312793752373
$ws.Range("H43").Value = 13.45737016776097
$ws.Range("G44").Value = 0.02438424464352845
$ws.Range("H44").Value = -31.67303636276195
$ws.Range("G45").Value = 0.007700669758706816
$ws.Range("H45").Value = -52.95914581098585
$ws.Range("G46").Value = 0.05991213956845752
$ws.Range("H46").Value = 65.32722910562802
$ws.Range("G47").Value = 0.08146847173690272
$ws.Range("H47").Value = 61.51460110398956
$ws.Range("G48").Value = 0.04482967537714999
$ws.Range("H48").Value = 4.788010289000823
$ws.Range("G49").Value = 0.04394657357260019
$ws.Range("H49").Value = -36.74576847551215
$ws.Range("G50").Value = 0.02730252103732629
$ws.Range("H50").Value = 58.06741726830158
$ws.Range("G51").Value = 0.02689722657693246
$ws.Range("H51").Value = 38.14907942277376
$ws.Range("G52").Value = -0.09698905763952953
$ws.Range("H52").Value = 6.309282302283532
$ws.Range("G53").Value = -0.06842636825123943
$ws.Range("H53").Value = 25.9094124426242
$ws.Range("G54").Value = 0.07778875156170577
$ws.Range("H54").Value = 6.381182601292282
$ws.Range("G55").Value = 0.09184749469313373
$ws.Range("H55").Value = 48.25671926391652
$ws.Range("G56").Value = 0.02209074913251056
$ws.Range("H56").Value = -36.86411892044853
$ws.Range("G57").Value = 0.0106872992631785
$ws.Range("H57").Value = 85.10897115566462
$ws.Range("G58").Value = 0.05433593940448143
$ws.Range("H58").Value = 117.2515825554734
$ws.Range("G59").Value = 0.02976558436578546
$ws.Range("H59").Value = 25.70655341974288
$ws.Range("G60").Value = 0.04101846235255185
$ws.Range("H60").Value = 26.43372918364897
$ws.Range("G61").Value = 0.03200459390181736
$ws.Range("H61").Value = 152.8409793911934
$ws.Range("G62").Value = 0.07155275983816951
$ws.Range("H62").Value = 18.53819849981627
$ws.Range("G63").Value = 0.07151866390372073
$ws.Range("H63").Value = 119.4531321711361
$ws.Range("G64").Value = 0.01924323859880504
$ws.Range("H64").Value = -52.51652605997614
$ws.Range("G65").Value = 0.04470234402438304
$ws.Range("H65").Value = -20.26255413963732
$ws.Range("G66").Value = 0.1162138588239151
$ws.Range("H66").Value = 24.22101556569537
$ws.Range("G67").Value = 0.1139807233203216
$ws.Range("H67").Value = -1.269792606126597
$ws.Range("G68").Value = -0.04298043932356216
$ws.Range("H68").Value = -23.32810275545842
$ws.Range("G69").Value = 0.009084140702235623
$ws.Range("H69").Value = 142.8055964742153
$ws.Range("G70").Value = 0.07136933824531512
$ws.Range("H70").Value = -22.95744687977239
$ws.Range("G71").Value = 0.06751550598069112
$ws.Range("H71").Value = -25.97697950071262
$ws.Range("G72").Value = -0.05651715716889315
$ws.Range("H72").Value = -0.7776801847205586
$ws.Range("G73").Value = -0.1095670998551015
$ws.Range("H73").Value = -48.54020180687038
$ws.Range("G74").Value = 0.1008428982636864
$ws.Range("H74").Value = 0.8965683820481898
$ws.Range("G75").Value = 0.1319076018571863
$ws.Range("H75").Value = 35.42315626662127
$ws.Range("G76").Value = 0.03241523938044558
$ws.Range("H76").Value = 26.76762162662265
$ws.Range("G77").Value = 0.01501380445696309
$ws.Range("H77").Value = 6.412656702651737
$ws.Range("G78").Value = 0.09114030370196438
$ws.Range("H78").Value = 41.79309582109754
$ws.Range("G79").Value = 0.09894927701230079
$ws.Range("H79").Value = 28.98481482872196
$ws.Range("G80").Value = -0.1968259108035692
$ws.Range("H80").Value = -18.85019602371191
$ws.Range("G81").Value = -0.1365362466125887
$ws.Range("H81").Value = 35.01121547656149
$ws.Range("G82").Value = 0.1292639796065003
$ws.Range("H82").Value = 12.70204448486605
$ws.Range("G83").Value = 0.1948298694351747
$ws.Range("H83").Value = 9.465921729562689
$ws.Range("G84").Value = 0.04848844290290327
$ws.Range("H84").Value = 103.411178732628
$ws.Range("G85").Value = 0.1016572206468263
$ws.Range("H85").Value = 65.0924375756844
